$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two sample/data rows (rows 2 and 3), keeping only the header row ---
$ws.Range("A2:A3").EntireRow.Delete()

# --- Rewrite the header row (row 1) to the new set of columns (A:H) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Image-text relation"
$ws.Range("C1").Value = "checkbox"
$ws.Range("D1").Value = "Decision part"
$ws.Range("E1").Value = "Hatefulness scale"
$ws.Range("F1").Value = "Confidence score"
$ws.Range("G1").Value = "Discard"
$ws.Range("H1").Value = "Elapsed Time (s)"

# Give the three new header cells (F1:H1) the same bold / bordered / centered
# header style already used by A1:E1 (copy formatting only, keep the values
# that were just written above).
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore the workbook selection to A2:H4 ---
$ws.Range("A2:H4").Select()
